# R6_Legacy_compare.xlsx -- documentation update
#
# "RLcomp_valid" (sheet: Construct/Plot VPC, Stepwise, Plot residuals) had
# three R6-notation cells whose `$` separators were not escaped with a
# backslash in the source text. Re-enter them with the escaped form
# (`\$`) to match the rest of the workbook's convention (e.g.
# "PM_data\$new()" on the first sheet).
$wb = $excel.ActiveWorkbook
$wsValid = $wb.Worksheets.Item("RLcomp_valid")

# Apply in the same order the original author touched them so that the
# shared-string table grows in the matching sequence.
$wsValid.Range("B6").Value = 'PM_result\$step()'
$wsValid.Range("B5").Value = 'PM_valid\$plot()'
$wsValid.Range("B3").Value = 'PM_result\$op\$plot(resid = T,…)'

# Restore the selection to a single cell (was a 2-cell range) without
# leaving that sheet the active one.
$wsValid.Range("B4").Select() | Out-Null

$wsOther = $wb.Worksheets.Item("RLcomp_other")
$wsOther.Activate() | Out-Null
